# Updated cryptos list -- refresh Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.640.15"
$ws.Range("E2").Value = "  +3.13%  "
$ws.Range("D3").Value = "3.408.99"
$ws.Range("E3").Value = "  +2.05%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "562.98"
$ws.Range("E5").Value = "  +2.62%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "176.10"
$ws.Range("E6").Value = "  +2.26%  "
$ws.Range("D8").Value = "3.403.78"
$ws.Range("E8").Value = "  +2.18%  "
$ws.Range("E9").Value = "  -0.04%  "
$ws.Range("E10").Value = "  +14.25%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.634"
$ws.Range("E11").Value = "  +3.33%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "55.07"
$ws.Range("E12").Value = "  +2.68%  "
$ws.Range("E13").Value = "  +6.26%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.18"
$ws.Range("E14").Value = "  +2.71%  "
$ws.Range("D15").Value = "3.949.52"
$ws.Range("E15").Value = "  +1.83%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "18.35"
$ws.Range("E16").Value = "  +2.47%  "
$ws.Range("D17").Value = "3.404.07"
$ws.Range("E17").Value = "  +2.17%  "
$ws.Range("E18").Value = "  +1.88%  "
$ws.Range("D19").Value = "65.555.89"
$ws.Range("E19").Value = "  +2.92%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.91"
$ws.Range("E20").Value = "  +1.55%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.996"
$ws.Range("E21").Value = "  +2.31%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "472.10"
$ws.Range("E22").Value = "  +14.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.18"
$ws.Range("E23").Value = "  +20.67%  "
$ws.Range("E24").Value = "  +2.58%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "86.76"
$ws.Range("E25").Value = "  +4.39%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.48"
$ws.Range("E26").Value = "  -1.11%  "
$ws.Range("E27").Value = "  +2.99%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.90"
$ws.Range("E28").Value = "  +6.50%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.91"
$ws.Range("E29").Value = "  +3.67%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "31.18"
$ws.Range("E30").Value = "  +7.14%  "
$ws.Range("E31").Value = "  +4.93%  "
$ws.Range("E32").Value = "  +2.11%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "62.75"
$ws.Range("E33").Value = "  +9.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "580.32"
$ws.Range("E34").Value = "  +1.33%  "
$ws.Range("E35").Value = "  +2.33%  "
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("E37").Value = "  -4.43%  "
$ws.Range("E38").Value = "  +4.66%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "35.98"
$ws.Range("E39").Value = "  +2.27%  "
$ws.Range("D40").Value = "0.0₃0762"
$ws.Range("E40").Value = "  +3.35%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.375"
$ws.Range("E41").Value = "  +2.30%  "
$ws.Range("D42").Value = "3.094.59"
$ws.Range("E42").Value = "  -1.64%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.998"
$ws.Range("E43").Value = "  -0.25%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.87"
$ws.Range("E44").Value = "  +2.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0419"
$ws.Range("E45").Value = "  +4.09%  "
$ws.Range("E46").Value = "  +4.03%  "
$ws.Range("E47").Value = "  +6.18%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.17"
$ws.Range("E48").Value = "  -2.65%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.60"
$ws.Range("E49").Value = "  +0.33%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.38"
$ws.Range("E50").Value = "  +4.48%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "136.83"
$ws.Range("E51").Value = "  +3.15%  "
